$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 104.5
$ws.Range("I5").Value = 99.36364
$ws.Range("K5").Value = 99.36364
$ws.Range("M5").Value = 15.63636
$ws.Range("H19").Value = 2233.4167
$ws.Range("J19").Value = 2466.6667
$ws.Range("L19").Value = 2466.6667
$ws.Range("N19").Value = -2816.6667
$ws.Range("H33").Value = 181.80952
$ws.Range("I33").Value = 182.36842
$ws.Range("J33").Value = 176.5
$ws.Range("K33").Value = 182.36842
$ws.Range("L33").Value = 176.5
$ws.Range("M33").Value = 46.63158000000001
$ws.Range("N33").Value = -634.5
$ws.Range("H40").Value = 47621076
$ws.Range("J40").Value = 55557724
$ws.Range("L40").Value = 55557724
$ws.Range("N40").Value = -55558074
$ws.Range("H76").Value = 39288610
$ws.Range("I76").Value = 45836250
$ws.Range("J76").Value = 2762.5
$ws.Range("K76").Value = 45836250
$ws.Range("L76").Value = 2762.5
$ws.Range("M76").Value = -45835935
$ws.Range("N76").Value = -3392.5
$ws.Range("H79").Value = 39288610
$ws.Range("I79").Value = 45836250
$ws.Range("J79").Value = 2762.5
$ws.Range("K79").Value = 45836250
$ws.Range("L79").Value = 2762.5
$ws.Range("M79").Value = -45835158
$ws.Range("N79").Value = -4946.5
$ws.Range("H116").Value = 4714.154
$ws.Range("I116").Value = 4853.778
$ws.Range("K116").Value = 4853.778
$ws.Range("M116").Value = -1411.778
$ws.Range("H132").Value = 1168.2778
$ws.Range("I132").Value = 1245.4565
$ws.Range("J132").Value = 724.5
$ws.Range("K132").Value = 3736.3695
$ws.Range("L132").Value = 2173.5
$ws.Range("M132").Value = -1206.3695
$ws.Range("N132").Value = -7233.5
$ws.Range("H135").Value = 1032.3846
$ws.Range("I135").Value = 596.90625
$ws.Range("J135").Value = 3023.1428
$ws.Range("K135").Value = 5372.15625
$ws.Range("L135").Value = 27208.2852
$ws.Range("M135").Value = -2837.15625
$ws.Range("N135").Value = -32278.2852
$ws.Range("H137").Value = 929.75
$ws.Range("I137").Value = 883.7778
$ws.Range("J137").Value = 1067.6666
$ws.Range("K137").Value = 2651.3334
$ws.Range("L137").Value = 3202.9998
$ws.Range("M137").Value = -101.3334
$ws.Range("N137").Value = -8302.9998
$ws.Range("H138").Value = 1360.22
$ws.Range("I138").Value = 616.0417
$ws.Range("J138").Value = 3273.8215
$ws.Range("K138").Value = 1848.1251
$ws.Range("L138").Value = 9821.4645
$ws.Range("M138").Value = 3291.8749
$ws.Range("N138").Value = -20101.4645

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18716.627
$ws.Range("I32").Value = 18175.143
$ws.Range("J32").Value = 20611.818
$ws.Range("K32").Value = 18175.143
$ws.Range("L32").Value = 20611.818
$ws.Range("M32").Value = -17888.143
$ws.Range("N32").Value = -21185.818
$ws.Range("H61").Value = 1151.5
$ws.Range("I61").Value = 719.34784
$ws.Range("K61").Value = 719.34784
$ws.Range("M61").Value = -507.34784
$ws.Range("H110").Value = 639.2273
$ws.Range("I110").Value = 570.7059
$ws.Range("J110").Value = 872.2
$ws.Range("K110").Value = 570.7059
$ws.Range("L110").Value = 872.2
$ws.Range("M110").Value = 1474.2941
$ws.Range("N110").Value = -4962.2
$ws.Range("H122").Value = 658.9375
$ws.Range("I122").Value = 580.2308
$ws.Range("K122").Value = 1740.6924
$ws.Range("M122").Value = 709.3075999999999
$ws.Range("H132").Value = 1510.9474
$ws.Range("I132").Value = 1227
$ws.Range("J132").Value = 1766.5
$ws.Range("K132").Value = 3681
$ws.Range("L132").Value = 5299.5
$ws.Range("M132").Value = -1151
$ws.Range("N132").Value = -10359.5
$ws.Range("H136").Value = 1151.5
$ws.Range("I136").Value = 719.34784
$ws.Range("K136").Value = 2158.04352
$ws.Range("M136").Value = 391.9564799999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 83732.60000000001
$ws.Range("I134").Value = 3959.7646
$ws.Range("J134").Value = 253249.88
$ws.Range("K134").Value = 11879.2938
$ws.Range("L134").Value = 759749.64
$ws.Range("M134").Value = -9344.293799999999
$ws.Range("N134").Value = -764819.64
$ws.Range("H140").Value = 95840
$ws.Range("J140").Value = 95840
$ws.Range("L140").Value = 95840
$ws.Range("N140").Value = -106200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H31").Value = 2720.2327
$ws.Range("I31").Value = 2677.9092
$ws.Range("J31").Value = 2859.9
$ws.Range("K31").Value = 2677.9092
$ws.Range("L31").Value = 2859.9
$ws.Range("M31").Value = -2382.9092
$ws.Range("N31").Value = -3449.9
$ws.Range("H34").Value = 2720.2327
$ws.Range("I34").Value = 2677.9092
$ws.Range("J34").Value = 2859.9
$ws.Range("K34").Value = 2677.9092
$ws.Range("L34").Value = 2859.9
$ws.Range("M34").Value = -2475.9092
$ws.Range("N34").Value = -3263.9
$ws.Range("H86").Value = 4650
$ws.Range("J86").Value = 4660
$ws.Range("L86").Value = 4660
$ws.Range("N86").Value = -6906
$ws.Range("H89").Value = 4650
$ws.Range("J89").Value = 4660
$ws.Range("L89").Value = 23300
$ws.Range("N89").Value = -34532
$ws.Range("H106").Value = 39700
$ws.Range("J106").Value = 39700
$ws.Range("L106").Value = 39700
$ws.Range("N106").Value = -42224
$ws.Range("H132").Value = 944
$ws.Range("I132").Value = 739.62195
$ws.Range("J132").Value = 2619.9
$ws.Range("K132").Value = 2218.86585
$ws.Range("L132").Value = 7859.700000000001
$ws.Range("M132").Value = 311.1341499999999
$ws.Range("N132").Value = -12919.7
$ws.Range("H134").Value = 1195.4814
$ws.Range("I134").Value = 997.97675
$ws.Range("J134").Value = 1967.5454
$ws.Range("K134").Value = 2993.93025
$ws.Range("L134").Value = 5902.6362
$ws.Range("M134").Value = -458.9302500000003
$ws.Range("N134").Value = -10972.6362
$ws.Range("H140").Value = 63704.285
$ws.Range("J140").Value = 63704.285
$ws.Range("L140").Value = 63704.285
$ws.Range("N140").Value = -74064.285

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 2000
$ws.Range("J19").Value = 2000
$ws.Range("L19").Value = 6000
$ws.Range("N19").Value = -6348
$ws.Range("H87").Value = 22450
$ws.Range("J87").Value = 23333.334
$ws.Range("L87").Value = 70000.00199999999
$ws.Range("N87").Value = -72496.00199999999
$ws.Range("H90").Value = 22450
$ws.Range("J90").Value = 23333.334
$ws.Range("L90").Value = 210000.006
$ws.Range("N90").Value = -222480.006
$ws.Range("H107").Value = 519102.34
$ws.Range("J107").Value = 778077.2
$ws.Range("L107").Value = 2334231.6
$ws.Range("N107").Value = -2338071.6
$ws.Range("H141").Value = 6140.2104
$ws.Range("I141").Value = 4208.769
$ws.Range("J141").Value = 10325
$ws.Range("K141").Value = 12626.307
$ws.Range("L141").Value = 30975
$ws.Range("M141").Value = -7446.307000000001
$ws.Range("N141").Value = -41335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 14961.4
$ws.Range("J17").Value = 14961.4
$ws.Range("L17").Value = 14961.4
$ws.Range("N17").Value = -15297.4
$ws.Range("H70").Value = 4093.7
$ws.Range("I70").Value = 3933.2144
$ws.Range("J70").Value = 4468.1665
$ws.Range("K70").Value = 3933.2144
$ws.Range("L70").Value = 4468.1665
$ws.Range("M70").Value = -3663.2144
$ws.Range("N70").Value = -5008.1665
$ws.Range("H73").Value = 4093.7
$ws.Range("I73").Value = 3933.2144
$ws.Range("J73").Value = 4468.1665
$ws.Range("K73").Value = 3933.2144
$ws.Range("L73").Value = 4468.1665
$ws.Range("M73").Value = -2997.2144
$ws.Range("N73").Value = -6340.1665
$ws.Range("H132").Value = 3180
$ws.Range("I132").Value = 3130
$ws.Range("J132").Value = 3300
$ws.Range("K132").Value = 9390
$ws.Range("L132").Value = 9900
$ws.Range("M132").Value = -6860
$ws.Range("N132").Value = -14960

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1545
$ws.Range("J46").Value = 1300
$ws.Range("L46").Value = 1300
$ws.Range("N46").Value = -1676
$ws.Range("H132").Value = 2033.4193
$ws.Range("I132").Value = 1889.8462
$ws.Range("K132").Value = 5669.5386
$ws.Range("M132").Value = -3139.5386
$ws.Range("H136").Value = 1510.1212
$ws.Range("I136").Value = 1447.4615
$ws.Range("J136").Value = 1742.8572
$ws.Range("K136").Value = 4342.3845
$ws.Range("L136").Value = 5228.571599999999
$ws.Range("M136").Value = -1792.3845
$ws.Range("N136").Value = -10328.5716
$ws.Range("H140").Value = 125629.336
$ws.Range("J140").Value = 125629.336
$ws.Range("L140").Value = 125629.336
$ws.Range("N140").Value = -135989.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 530.5
$ws.Range("I113").Value = 526
$ws.Range("J113").Value = 538.6
$ws.Range("K113").Value = 1578
$ws.Range("L113").Value = 1615.8
$ws.Range("M113").Value = 592
$ws.Range("N113").Value = -5955.8
$ws.Range("H132").Value = 1404.439
$ws.Range("I132").Value = 895.44446
$ws.Range("J132").Value = 2386.0715
$ws.Range("K132").Value = 2686.33338
$ws.Range("L132").Value = 7158.2145
$ws.Range("M132").Value = -156.33338
$ws.Range("N132").Value = -12218.2145
$ws.Range("H136").Value = 656.3889
$ws.Range("I136").Value = 332.1875
$ws.Range("J136").Value = 3250
$ws.Range("K136").Value = 3250
$ws.Range("L136").Value = 9750
$ws.Range("M136").Value = 1553.4375
$ws.Range("N136").Value = -14850
